$d = $word.ActiveDocument

# Scope edits to the paragraph that lists the three research thrusts
# ("We propose to develop 1) ... 2) ... 3) ...") to avoid touching any
# other similarly-worded text elsewhere in the document.
$p = $d.Paragraphs(7)
$r = $p.Range

# --- Change 1 -----------------------------------------------------------
# "...computer architects, researchers, and software developers..."
#   -> "...computer architects and software developers..."
$r.Find.Execute(", researchers, and software developers", $true, $false, $false, $false, $false, $true, 1, $false, " and software developers", 2)

# --- Change 2 -----------------------------------------------------------
# Rewrite research item 2) and item 3) to match the updated project
# description.
$old2 = "2) methods for efficient circuit-level exploration of caches and functional units that can be integrated into architecture-level simulators, analogous to how Cacti and McPat are used to obtain per-event latency and power estimates in cycle-accurate simulators, and 3) methods for " + [char]0x201C + "calibration" + [char]0x201D + " of simulation parameters against measured signals"
$new2 = "2) methods for efficient circuit-level simulation of microarchitectural building blocks to generate per-event current/voltage shapes, 3) methods for transforming and scaling current/voltage shapes profiles into per-event signal-snippets that will be " + [char]0x201C + "stitched" + [char]0x201D + " into longer signals through cycle-accurate simulation, and 3) methods for calibration of the simulation against measured signals"

$r = $p.Range
$r.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
